$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Absent"
$ws.Range("C5").Value = "Absent"
$ws.Range("C6").Value = "Present"
$ws.Range("C7").Value = "Present"
$ws.Range("C8").Value = "Present"
$ws.Range("C9").Value = "Present"
$ws.Range("C10").Value = "Present"
